$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell's text (Price in column D, Volume(1h) % in column E).
# NumberFormat is forced to Text ("@") before assignment so numeric-looking strings
# (e.g. "28.167.68", "338.53") are stored as text, matching the source data which
# keeps these as plain text cells rather than numbers.
$changes = @(
    @{ Cell = 'D2'; Value = '28.167.68' }
    @{ Cell = 'E2'; Value = '  +0.11%  ' }
    @{ Cell = 'D3'; Value = '1.799.18' }
    @{ Cell = 'E3'; Value = '  +2.19%  ' }
    @{ Cell = 'E4'; Value = '  +0.28%  ' }
    @{ Cell = 'D5'; Value = '338.53' }
    @{ Cell = 'E5'; Value = '  +1.12%  ' }
    @{ Cell = 'D6'; Value = '1.002' }
    @{ Cell = 'E6'; Value = '  +0.40%  ' }
    @{ Cell = 'D7'; Value = '0.4658' }
    @{ Cell = 'E7'; Value = '  +23.00%  ' }
    @{ Cell = 'D8'; Value = '0.3723' }
    @{ Cell = 'E8'; Value = '  +10.61%  ' }
    @{ Cell = 'D9'; Value = '45.58' }
    @{ Cell = 'D10'; Value = '0.07691' }
    @{ Cell = 'E10'; Value = '  +6.62%  ' }
    @{ Cell = 'D11'; Value = '1.147' }
    @{ Cell = 'E11'; Value = '  +1.59%  ' }
    @{ Cell = 'D12'; Value = '22.62' }
    @{ Cell = 'E12'; Value = '  -0.30%  ' }
    @{ Cell = 'D13'; Value = '1.004' }
    @{ Cell = 'E13'; Value = '  +0.36%  ' }
    @{ Cell = 'D14'; Value = '6.408' }
    @{ Cell = 'E14'; Value = '  +3.59%  ' }
    @{ Cell = 'D15'; Value = '7.448' }
    @{ Cell = 'E15'; Value = '  +3.06%  ' }
    @{ Cell = 'D16'; Value = '1.797.98' }
    @{ Cell = 'E16'; Value = '  +2.29%  ' }
    @{ Cell = 'D17'; Value = '0.00001097' }
    @{ Cell = 'E17'; Value = '  +3.84%  ' }
    @{ Cell = 'D18'; Value = '0.06762' }
    @{ Cell = 'E18'; Value = '  +2.82%  ' }
    @{ Cell = 'D19'; Value = '82.35' }
    @{ Cell = 'E19'; Value = '  +1.71%  ' }
    @{ Cell = 'E20'; Value = '  +0.35%  ' }
    @{ Cell = 'D21'; Value = '17.54' }
    @{ Cell = 'E21'; Value = '  +3.02%  ' }
    @{ Cell = 'D22'; Value = '6.451' }
    @{ Cell = 'E22'; Value = '  +2.77%  ' }
    @{ Cell = 'D23'; Value = '28.158.13' }
    @{ Cell = 'E23'; Value = '  +0.12%  ' }
    @{ Cell = 'D24'; Value = '11.98' }
    @{ Cell = 'E24'; Value = '  +2.76%  ' }
    @{ Cell = 'D25'; Value = '2.404' }
    @{ Cell = 'E25'; Value = '  +0.35%  ' }
    @{ Cell = 'D26'; Value = '20.95' }
    @{ Cell = 'E26'; Value = '  +5.19%  ' }
    @{ Cell = 'D27'; Value = '2.412' }
    @{ Cell = 'E27'; Value = '  +3.25%  ' }
    @{ Cell = 'D28'; Value = '151.51' }
    @{ Cell = 'E28'; Value = '  -0.97%  ' }
    @{ Cell = 'D29'; Value = '2.008.51' }
    @{ Cell = 'E29'; Value = '  +2.53%  ' }
    @{ Cell = 'D30'; Value = '134.84' }
    @{ Cell = 'E30'; Value = '  +1.75%  ' }
    @{ Cell = 'E31'; Value = '  +0.21%  ' }
    @{ Cell = 'D32'; Value = '4.061' }
    @{ Cell = 'E32'; Value = '  +1.08%  ' }
    @{ Cell = 'D33'; Value = '5.966' }
    @{ Cell = 'E33'; Value = '  +2.52%  ' }
    @{ Cell = 'D34'; Value = '0.09647' }
    @{ Cell = 'E34'; Value = '  +9.84%  ' }
    @{ Cell = 'D35'; Value = '0.02399' }
    @{ Cell = 'E35'; Value = '  +2.31%  ' }
    @{ Cell = 'D36'; Value = '0.2231' }
    @{ Cell = 'E36'; Value = '  +5.37%  ' }
    @{ Cell = 'D37'; Value = '12.22' }
    @{ Cell = 'E37'; Value = '  -0.68%  ' }
    @{ Cell = 'D38'; Value = '0.06406' }
    @{ Cell = 'E38'; Value = '  +3.04%  ' }
    @{ Cell = 'D39'; Value = '5.287' }
    @{ Cell = 'E39'; Value = '  +2.22%  ' }
    @{ Cell = 'E40'; Value = '  +1.16%  ' }
    @{ Cell = 'D41'; Value = '1.239' }
    @{ Cell = 'E41'; Value = '  +1.63%  ' }
    @{ Cell = 'E42'; Value = '  +2.52%  ' }
    @{ Cell = 'D43'; Value = '8.102' }
    @{ Cell = 'E43'; Value = '  +1.05%  ' }
    @{ Cell = 'D44'; Value = '14.34' }
    @{ Cell = 'E44'; Value = '  +4.74%  ' }
    @{ Cell = 'E45'; Value = '  +0.39%  ' }
    @{ Cell = 'D46'; Value = '0.6166' }
    @{ Cell = 'E46'; Value = '  +1.77%  ' }
    @{ Cell = 'D47'; Value = '3.869' }
    @{ Cell = 'E47'; Value = '  +1.47%  ' }
    @{ Cell = 'D48'; Value = '130.06' }
    @{ Cell = 'E48'; Value = '  +0.30%  ' }
    @{ Cell = 'D49'; Value = '2.075' }
    @{ Cell = 'E49'; Value = '  +3.08%  ' }
    @{ Cell = 'D50'; Value = '1.185' }
    @{ Cell = 'E50'; Value = '  -0.66%  ' }
    @{ Cell = 'D51'; Value = '0.07126' }
    @{ Cell = 'E51'; Value = '  -1.20%  ' }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
}
